# Update RF (column I) values for rows 36-83 on Sheet1 from 65.87444444444445
# to 36.83563909774436, reflecting the 2025 data / RF changes described in
# the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$newRF = 36.83563909774436

$ws.Range("I36:I83").Value = $newRF
